$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of the header row.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row (2-42)
$ws.Range("AD2:AD42").Value = 90
$ws.Range("AE2:AE42").Value = 72
$ws.Range("AF2:AF42").Value = 0
